$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the "last row" emphasis formatting (bold bottom border) from the
#    soon-to-be-deleted worker row (49) onto row 48, which will become the
#    new last row of the periods table once row 49 is removed.
# ---------------------------------------------------------------------------
$ws.Range("B49:J49").Copy() | Out-Null
$ws.Range("B48:J48").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Rewrite the periods table (rows 16-48) in ascending order 1607..1903.
#    The "Valor Mora" amounts stay attached to the table positions: every
#    row keeps 24640 except the final row, which carries 20533.
# ---------------------------------------------------------------------------
$periods = @(
    "1607","1608","1609","1610","1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    if ($row -eq 48) {
        $ws.Range("F$row").Value = 20533
    } else {
        $ws.Range("F$row").Value = 24640
    }
    $ws.Range("G$row").Value = 616000
}

# ---------------------------------------------------------------------------
# 3) Remove the second worker (ARNOLD YESITH SIERRA TORRES) row entirely.
# ---------------------------------------------------------------------------
$ws.Rows.Item(49).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4) Update the summary figures: total overdue value and worker count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 809013
$ws.Range("C13").Value = 1

# ---------------------------------------------------------------------------
# 5) Header row fix-up.
# ---------------------------------------------------------------------------
$ws.Range("I15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 6) Column D is now narrower since the remaining worker name is shorter.
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 27.2
